$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.331.19'
$ws.Range("E2").Value = '  +0.19%  '
$ws.Range("D3").Value = '1.668.98'
$ws.Range("E3").Value = '  +0.30%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '220.28'
$ws.Range("E5").Value = '  +0.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5282'
$ws.Range("E6").Value = '  -0.83%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2648'
$ws.Range("E8").Value = '  +0.21%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06369'
$ws.Range("E9").Value = '  +0.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.91'
$ws.Range("E10").Value = '  +1.66%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07849'
$ws.Range("E11").Value = '  +0.01%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.526'
$ws.Range("E12").Value = '  -0.54%  '
$ws.Range("D13").Value = '1.673.39'
$ws.Range("E13").Value = '  +1.93%  '
$ws.Range("D14").Value = '1.899.36'
$ws.Range("E14").Value = '  +0.36%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5610'
$ws.Range("E15").Value = '  +1.42%  '
$ws.Range("D16").Value = '0.0₅8109'
$ws.Range("E16").Value = '  -1.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.96'
$ws.Range("E17").Value = '  +0.37%  '
$ws.Range("D18").Value = '26.346.19'
$ws.Range("E18").Value = '  +0.23%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.011'
$ws.Range("E19").Value = '  +0.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.718'
$ws.Range("E20").Value = '  +0.91%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '200.36'
$ws.Range("E21").Value = '  +4.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.29'
$ws.Range("E22").Value = '  +0.73%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.055'
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("E24").Value = '  +0.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.70'
$ws.Range("E25").Value = '  +1.08%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1212'
$ws.Range("E26").Value = '  -1.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.251'
$ws.Range("E27").Value = '  -0.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.22'
$ws.Range("E28").Value = '  +0.52%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.512'
$ws.Range("E29").Value = '  +2.72%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05886'
$ws.Range("E30").Value = '  +0.61%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.283'
$ws.Range("E31").Value = '  +0.37%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.511'
$ws.Range("E32").Value = '  -2.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.323'
$ws.Range("E33").Value = '  +0.37%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.599'
$ws.Range("E34").Value = '  -1.30%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9645'
$ws.Range("E35").Value = '  +0.52%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.823'
$ws.Range("E36").Value = '  +0.07%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.436'
$ws.Range("E37").Value = '  +0.27%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5792'
$ws.Range("E38").Value = '  -0.48%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01617'
$ws.Range("E39").Value = '  +0.07%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.950'
$ws.Range("E40").Value = '  +1.02%  '
$ws.Range("D41").Value = '1.076.46'
$ws.Range("E41").Value = '  +2.71%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8596'
$ws.Range("E42").Value = '  +0.55%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.011'
$ws.Range("E43").Value = '  +0.14%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '102.97'
$ws.Range("E44").Value = '  -1.51%  '
$ws.Range("D45").Value = '1.810.11'
$ws.Range("E45").Value = '  +0.25%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '58.50'
$ws.Range("E46").Value = '  +1.92%  '
$ws.Range("E47").Value = '  -1.18%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.014'
$ws.Range("E48").Value = '  +0.47%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4421'
$ws.Range("E49").Value = '  +1.11%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.029'
$ws.Range("E50").Value = '  +0.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05155'
$ws.Range("E51").Value = '  -0.16%  '
